# Update bus voltage magnitude results (vm_pu) for the 380 kV case.
# Slack/reference bus voltage (column B) drops from 1.05 p.u. to 1.02 p.u.,
# and the resulting steady-state voltages for the other buses are updated accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.011088528027686
$ws.Range("D2").Value = 1.039542189857305
$ws.Range("E2").Value = 1.013307971452307
$ws.Range("F2").Value = 1.041759721710385
$ws.Range("I2").Value = 1.033135816629002
$ws.Range("J2").Value = 1.016339018249044
$ws.Range("K2").Value = 1.042327081931044
$ws.Range("L2").Value = 1.016169514477233
$ws.Range("M2").Value = 1.044538326767334
$ws.Range("N2").Value = 1.00956726760115
# row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.011998492306645
$ws.Range("D3").Value = 1.04008428611555
$ws.Range("E3").Value = 1.014078504524762
$ws.Range("F3").Value = 1.042529858890228
$ws.Range("I3").Value = 1.033199422137062
$ws.Range("J3").Value = 1.016882295147649
$ws.Range("K3").Value = 1.042680125176134
$ws.Range("L3").Value = 1.016744843737334
$ws.Range("M3").Value = 1.045119263673284
$ws.Range("N3").Value = 1.009750598916128
# row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.012588144329397
$ws.Range("D4").Value = 1.040434938040297
$ws.Range("E4").Value = 1.01457820428877
$ws.Range("F4").Value = 1.043028640504153
$ws.Range("I4").Value = 1.033239147806063
$ws.Range("J4").Value = 1.017234091721952
$ws.Range("K4").Value = 1.042907738730191
$ws.Range("L4").Value = 1.017117577590346
$ws.Range("M4").Value = 1.045494959480256
$ws.Range("N4").Value = 1.009869210989908
# row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.012836234969218
$ws.Range("D5").Value = 1.040582321047545
$ws.Range("E5").Value = 1.014788543049754
$ws.Range("F5").Value = 1.043238434032102
$ws.Range("I5").Value = 1.033255505014807
$ws.Range("J5").Value = 1.017382048193519
$ws.Range("K5").Value = 1.043003227337204
$ws.Range("L5").Value = 1.017274383212267
$ws.Range("M5").Value = 1.045652849923052
$ws.Range("N5").Value = 1.009919071338718
# row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.012877902227306
$ws.Range("D6").Value = 1.040607065417681
$ws.Range("E6").Value = 1.014823875346598
$ws.Range("F6").Value = 1.043273665421547
$ws.Range("I6").Value = 1.033258231298179
$ws.Range("J6").Value = 1.017406894293229
$ws.Range("K6").Value = 1.043019248510365
$ws.Range("L6").Value = 1.017300717883226
$ws.Range("M6").Value = 1.045679357311604
$ws.Range("N6").Value = 1.009927442853383
# row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.012591458540636
$ws.Range("D7").Value = 1.040436907504079
$ws.Range("E7").Value = 1.01458101380635
$ws.Range("F7").Value = 1.043031443362998
$ws.Range("I7").Value = 1.033239367722461
$ws.Range("J7").Value = 1.01723606848447
$ws.Range("K7").Value = 1.04290901544211
$ws.Range("L7").Value = 1.017119672410636
$ws.Range("M7").Value = 1.045497069427211
$ws.Range("N7").Value = 1.009869877243261
# row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.011395879180958
$ws.Range("D8").Value = 1.039725417483242
$ws.Range("E8").Value = 1.013568145251036
$ws.Range("F8").Value = 1.042019898325313
$ws.Range("I8").Value = 1.033157608300696
$ws.Range("J8").Value = 1.016522566719628
$ws.Range("K8").Value = 1.042446565209147
$ws.Range("L8").Value = 1.016363854286721
$ws.Range("M8").Value = 1.044734699240397
$ws.Range("N8").Value = 1.00962922822854
# row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.009295634788565
$ws.Range("D9").Value = 1.038470854025741
$ws.Range("E9").Value = 1.01179193715296
$ws.Range("F9").Value = 1.040240986010304
$ws.Range("I9").Value = 1.033002610410435
$ws.Range("J9").Value = 1.015267325020622
$ws.Range("K9").Value = 1.041625392571447
$ws.Range("L9").Value = 1.015035561170542
$ws.Range("M9").Value = 1.043389790929709
$ws.Range("N9").Value = 1.009205074539439
# row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.007899925812085
$ws.Range("D10").Value = 1.037634052663084
$ws.Range("E10").Value = 1.010613670428612
$ws.Range("F10").Value = 1.039057575734623
$ws.Range("I10").Value = 1.032891978863282
$ws.Range("J10").Value = 1.014431932455074
$ws.Range("K10").Value = 1.041073823101172
$ws.Range("L10").Value = 1.014152492881109
$ws.Range("M10").Value = 1.04249228540362
$ws.Range("N10").Value = 1.008922265092357
# row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.00729663956549
$ws.Range("D11").Value = 1.037271635049517
$ws.Range("E11").Value = 1.010104880235068
$ws.Range("F11").Value = 1.038545775612385
$ws.Range("I11").Value = 1.032842352779805
$ws.Range("J11").Value = 1.014070552087649
$ws.Range("K11").Value = 1.04083403083487
$ws.Range("L11").Value = 1.013770712929152
$ws.Range("M11").Value = 1.042103467139973
$ws.Range("N11").Value = 1.008799801487448
# row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.007072713350215
$ws.Range("D12").Value = 1.037137007664684
$ws.Range("E12").Value = 1.009916105880841
$ws.Range("F12").Value = 1.03835576642473
$ws.Range("I12").Value = 1.032823661535442
$ws.Range("J12").Value = 1.013936373002746
$ws.Range("K12").Value = 1.040744818903479
$ws.Range("L12").Value = 1.013628993351439
$ws.Range("M12").Value = 1.041959015829301
$ws.Range("N12").Value = 1.00875431270911
# row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.007120738972072
$ws.Range("D13").Value = 1.037165886082851
$ws.Range("E13").Value = 1.009956588957568
$ws.Range("F13").Value = 1.038396519650847
$ws.Range("I13").Value = 1.032827682529029
$ws.Range("J13").Value = 1.013965152425855
$ws.Range("K13").Value = 1.040763961582915
$ws.Range("L13").Value = 1.013659388567776
$ws.Range("M13").Value = 1.041990002310518
$ws.Range("N13").Value = 1.008764070214834
# row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.007278126460365
$ws.Range("D14").Value = 1.037260506885649
$ws.Range("E14").Value = 1.010089271734238
$ws.Range("F14").Value = 1.038530067418245
$ws.Range("I14").Value = 1.032840813012398
$ws.Range("J14").Value = 1.014059459701496
$ws.Range("K14").Value = 1.040826659444162
$ws.Range("L14").Value = 1.013758996491707
$ws.Range("M14").Value = 1.042091527286578
$ws.Range("N14").Value = 1.008796041375255
# row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.007375119476205
$ws.Range("D15").Value = 1.037318804713395
$ws.Range("E15").Value = 1.010171050240908
$ws.Range("F15").Value = 1.038612363418935
$ws.Range("I15").Value = 1.032848868984866
$ws.Range("J15").Value = 1.014117572670962
$ws.Range("K15").Value = 1.040865270845583
$ws.Range("L15").Value = 1.013820380248122
$ws.Range("M15").Value = 1.04215407667144
$ws.Range("N15").Value = 1.008815739830283
# row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.007939986469644
$ws.Range("D16").Value = 1.03765810369322
$ws.Range("E16").Value = 1.010647466938433
$ws.Range("F16").Value = 1.039091555579309
$ws.Range("I16").Value = 1.032895236171959
$ws.Range("J16").Value = 1.014455923548062
$ws.Range("K16").Value = 1.041089717261592
$ws.Range("L16").Value = 1.014177842984687
$ws.Range("M16").Value = 1.042518086108919
$ws.Range("N16").Value = 1.008930392525961
# row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.008294598546172
$ws.Range("D17").Value = 1.037870917918794
$ws.Range("E17").Value = 1.010946688375274
$ws.Range("F17").Value = 1.039392308986026
$ws.Range("I17").Value = 1.032923860640211
$ws.Range("J17").Value = 1.014668256608771
$ws.Range("K17").Value = 1.0412302509917
$ws.Range("L17").Value = 1.014402229902357
$ws.Range("M17").Value = 1.042746369514125
$ws.Range("N17").Value = 1.009002310042415
# row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.008501540444878
$ws.Range("D18").Value = 1.037995041267277
$ws.Range("E18").Value = 1.011121354833115
$ws.Range("F18").Value = 1.039567793440808
$ws.Range("I18").Value = 1.032940390617763
$ws.Range("J18").Value = 1.014792140600834
$ws.Range("K18").Value = 1.041312129292504
$ws.Range("L18").Value = 1.014533168270762
$ws.Range("M18").Value = 1.042879504658586
$ws.Range("N18").Value = 1.009044257776758
# row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.008572119664547
$ws.Range("D19").Value = 1.038037362766112
$ws.Range("E19").Value = 1.011180934543009
$ws.Range("F19").Value = 1.039627639223689
$ws.Range("I19").Value = 1.032945998700158
$ws.Range("J19").Value = 1.014834387525931
$ws.Range("K19").Value = 1.041340031911579
$ws.Range("L19").Value = 1.0145778245495
$ws.Range("M19").Value = 1.042924897068434
$ws.Range("N19").Value = 1.009058560765627
# row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.008256541397925
$ws.Range("D20").Value = 1.037848085735229
$ws.Range("E20").Value = 1.010914570703941
$ws.Range("F20").Value = 1.039360034760243
$ws.Range("I20").Value = 1.032920806689842
$ws.Range("J20").Value = 1.014645471786319
$ws.Range("K20").Value = 1.041215182623174
$ws.Range("L20").Value = 1.014378149376809
$ws.Range("M20").Value = 1.042721878785221
$ws.Range("N20").Value = 1.008994594025769
# row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.007231775309732
$ws.Range("D21").Value = 1.037232643657673
$ws.Range("E21").Value = 1.010050194071622
$ws.Range("F21").Value = 1.038490738255719
$ws.Range("I21").Value = 1.032836953525503
$ws.Range("J21").Value = 1.014031687069087
$ws.Range("K21").Value = 1.040808200404402
$ws.Range("L21").Value = 1.013729661933126
$ws.Range("M21").Value = 1.042061631437837
$ws.Range("N21").Value = 1.008786626671255
# row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.006588397187797
$ws.Range("D22").Value = 1.036845637821808
$ws.Range("E22").Value = 1.009507959180963
$ws.Range("F22").Value = 1.037944734593239
$ws.Range("I22").Value = 1.032782739929191
$ws.Range("J22").Value = 1.013646087329338
$ws.Range("K22").Value = 1.040551491434683
$ws.Range("L22").Value = 1.013322456304372
$ws.Range("M22").Value = 1.04164635311206
$ws.Range("N22").Value = 1.008655867578719
# row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.00692937527705
$ws.Range("D23").Value = 1.037050801214461
$ws.Range("E23").Value = 1.009795290753383
$ws.Range("F23").Value = 1.038234127827253
$ws.Range("I23").Value = 1.032811620723649
$ws.Range("J23").Value = 1.013850471202789
$ws.Range("K23").Value = 1.040687655177495
$ws.Range("L23").Value = 1.013538273631423
$ws.Range("M23").Value = 1.041866513992974
$ws.Range("N23").Value = 1.008725185465943
# row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.008273737467932
$ws.Range("D24").Value = 1.037858402639911
$ws.Range("E24").Value = 1.010929082878165
$ws.Range("F24").Value = 1.039374617908294
$ws.Range("I24").Value = 1.032922187152152
$ws.Range("J24").Value = 1.014655767162784
$ws.Range("K24").Value = 1.041221991656764
$ws.Range("L24").Value = 1.014389030152935
$ws.Range("M24").Value = 1.042732945149453
$ws.Range("N24").Value = 1.008998080563269
# row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.009837818476399
$ws.Range("D25").Value = 1.038795273305037
$ws.Range("E25").Value = 1.012250101332534
$ws.Range("F25").Value = 1.040700441431569
$ws.Range("I25").Value = 1.033043971146071
$ws.Range("J25").Value = 1.015591586805981
$ws.Range("K25").Value = 1.041838418953857
$ws.Range("L25").Value = 1.01537852825826
$ws.Range("M25").Value = 1.043737649243696
$ws.Range("N25").Value = 1.009314737482208
